$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 32.47042033333333
$ws.Range("H2").Value = 97.411261
$ws.Range("I2").Value = 0.5240295449207956
$ws.Range("J2").Value = 0.5240295449207955
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 45.95651366666667
$ws.Range("N2").Value = 137.869541
$ws.Range("O2").Value = 0.6189188856627118
$ws.Range("P2").Value = 0.6189188856627118
$ws.Range("Q2").Value = 1492.227315811244
$ws.Range("R2").Value = 13430.0458423012
$ws.Range("S2").Value = 0.3243317819967168
$ws.Range("T2").Value = 0.3243317819967167

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 32.47042033333333
$ws.Range("H3").Value = 97.411261
$ws.Range("I3").Value = 0.5240295449207956
$ws.Range("J3").Value = 0.5240295449207955
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("N3").Value = 20.549744
$ws.Range("O3").Value = 0.09225115688993263
$ws.Range("P3").Value = 0.09225115688993261
$ws.Range("Q3").Value = 222.4196084741315
$ws.Range("R3").Value = 2001.776476267184
$ws.Range("S3").Value = 0.04834233176344831
$ws.Range("T3").Value = 0.04834233176344829

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 32.47042033333333
$ws.Range("H4").Value = 97.411261
$ws.Range("I4").Value = 0.5240295449207956
$ws.Range("J4").Value = 0.5240295449207955
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 21.446458
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2888299574473556
$ws.Range("P4").Value = 0.2888299574473556
$ws.Range("Q4").Value = 696.3755059211792
$ws.Range("R4").Value = 6267.379553290613
$ws.Range("S4").Value = 0.1513554311606305
$ws.Range("T4").Value = 0.1513554311606305

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 22.43791
$ws.Range("H5").Value = 67.31372999999999
$ws.Range("I5").Value = 0.3621181261458191
$ws.Range("J5").Value = 0.362118126145819
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 45.95651366666667
$ws.Range("N5").Value = 137.869541
$ws.Range("O5").Value = 0.6189188856627118
$ws.Range("P5").Value = 0.6189188856627118
$ws.Range("Q5").Value = 1031.168117566437
$ws.Range("R5").Value = 9280.513058097929
$ws.Range("S5").Value = 0.2241217471124397
$ws.Range("T5").Value = 0.2241217471124396

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 22.43791
$ws.Range("H6").Value = 67.31372999999999
$ws.Range("I6").Value = 0.3621181261458191
$ws.Range("J6").Value = 0.362118126145819
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("N6").Value = 20.549744
$ws.Range("O6").Value = 0.09225115688993263
$ws.Range("P6").Value = 0.09225115688993261
$ws.Range("Q6").Value = 153.6977687983467
$ws.Range("R6").Value = 1383.27991918512
$ws.Range("S6").Value = 0.03340581606776637
$ws.Range("T6").Value = 0.03340581606776635

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 22.43791
$ws.Range("H7").Value = 67.31372999999999
$ws.Range("I7").Value = 0.3621181261458191
$ws.Range("J7").Value = 0.362118126145819
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 21.446458
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2888299574473556
$ws.Range("P7").Value = 0.2888299574473556
$ws.Range("Q7").Value = 481.2136944227799
$ws.Range("R7").Value = 4330.923249805019
$ws.Range("S7").Value = 0.1045905629656131
$ws.Range("T7").Value = 0.104590562965613

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 7.054627
$ws.Range("H8").Value = 21.163881
$ws.Range("I8").Value = 0.1138523289333856
$ws.Range("J8").Value = 0.1138523289333855
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 45.95651366666667
$ws.Range("N8").Value = 137.869541
$ws.Range("O8").Value = 0.6189188856627118
$ws.Range("P8").Value = 0.6189188856627118
$ws.Range("Q8").Value = 324.2060621387357
$ws.Range("R8").Value = 2917.854559248621
$ws.Range("S8").Value = 0.07046535655355551
$ws.Range("T8").Value = 0.0704653565535555

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 7.054627
$ws.Range("H9").Value = 21.163881
$ws.Range("I9").Value = 0.1138523289333856
$ws.Range("J9").Value = 0.1138523289333855
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("N9").Value = 20.549744
$ws.Range("O9").Value = 0.09225115688993263
$ws.Range("P9").Value = 0.09225115688993261
$ws.Range("Q9").Value = 48.32359295516267
$ws.Range("R9").Value = 434.912336596464
$ws.Range("S9").Value = 0.01050300905871797
$ws.Range("T9").Value = 0.01050300905871796

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 7.054627
$ws.Range("H10").Value = 21.163881
$ws.Range("I10").Value = 0.1138523289333856
$ws.Range("J10").Value = 0.1138523289333855
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 21.446458
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2888299574473556
$ws.Range("P10").Value = 0.2888299574473556
$ws.Range("Q10").Value = 151.296761661166
$ws.Range("R10").Value = 1361.670854950494
$ws.Range("S10").Value = 0.03288396332111208
$ws.Range("T10").Value = 0.03288396332111208

Write-Output "Update complete"
